$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 14: "suporte pilha 9v" component -----------------------------
# (E14 is populated first so its URL text lands in the shared-strings table
# ahead of the plain "suporte pilha 9v" label, matching the order the
# strings were authored in.)
$ws.Range("E14").Value = "https://www.ptrobotics.com/suporte-pilhas/366-pp3-connector.html?gclid=Cj0KCQjwguGYBhDRARIsAHgRm4-VrK0A46wuEH-PEt_CZW79EWvXwBWDFZ_1E3wm6hctYjKeHSZoDbIaAte6EALw_wcB"
$ws.Range("B14").Value = "suporte pilha 9v"
$ws.Range("C14").Value = 1

# Copy the number-format/style used by the other "Custo por unidade" cells
# (D7:D13) onto D14, then set its value -- keeps cellXfs index 1 (style
# "Moeda", numFmt 164) instead of minting a brand new style entry.
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Value = 0.43

# --- Extend the total formula to include the new row -----------------------
$ws.Range("D15").Formula = "=SUM(D7:D14)"

# --- Turn the Transistor's existing URL text into a real hyperlink --------
$ws.Hyperlinks.Add($ws.Range("E10"), "https://eu.mouser.com/ProductDetail/onsemi-Fairchild/BC547CBU?qs=SpPkH8nd0tbX52Pj0mX5Aw%3D%3D") | Out-Null

# --- Underline formatting on E17 (matches the style picked up next to the
#     other totals-area cells) -----------------------------------------------
$ws.Range("E17").Font.Underline = $true

# --- View: zoom the sheet to 115% ------------------------------------------
$ws.Select()
$excel.ActiveWindow.Zoom = 115
